# ProfileDatabase.xlsx edit
#
# Rows 2 and 3 of the "profiles" sheet effectively swapped places:
#   - Row 2 now holds the data that used to be in row 3 (RichDogeyBoy / upland999 / ...)
#   - Row 3 now holds the data that used to be in row 2 (trashboatsr / dogeyboy19 / ...),
#     but with a refreshed Balance (D) and a refreshed Bearer Token (E) - a new JWT was
#     issued (new tokenId + iat) - and the trailing LichessReplaced (I) value was not
#     carried over / left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: becomes RichDogeyBoy's profile (previously row 3's content), values only,
#     keep the row's existing formatting untouched. ---
$ws.Range("A2").Value = "RichDogeyBoy"
$ws.Range("B2").Value = "upland999"
$ws.Range("C2").Value = 1500
$ws.Range("D2").Value = 806509
$ws.Range("E2").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VySWQiOiI3NTUyYjg4MC03MGY0LTExZWQtYWQ1YS0yN2QxMzZmYWRhMGUiLCJhcHBJZCI6MjMyLCJ0b2tlbklkIjoiNzNkNTk1ZmUtYjFiYy00MGY5LWJiMzctN2VhNGU0MTk1NTIyIiwiaWF0IjoxNzE1NTgwNjgyfQ.m1DYoR_QvwWGfMGoAKi6f_H0nCN7U_zanZ9giuWGEys"
$ws.Range("F2").Value = "kzolrwakubj3"
$ws.Range("G2").Value = "carlos"
$ws.Range("H2").Value = "7552b880-70f4-11ed-ad5a-27d136fada0e"
$ws.Range("I2").Value = 4

# --- Row 3: drop the old row entirely (its content, style and the trailing
#     LichessReplaced cell all go away), then refill it with trashboatsr's
#     profile, whose Balance and Bearer Token have changed. ---
$ws.Rows("3:3").Delete()
$ws.Rows("3:3").RowHeight = 19.5

$ws.Range("A3").Value = "trashboatsr"
$ws.Range("B3").Value = "dogeyboy19"
$ws.Range("C3").Value = 1890
$ws.Range("D3").Value = 1891734
$ws.Range("E3").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VySWQiOiI0ODI5OGVhMC0yNDBhLTExZWUtOWMwNC1iMzcyMDk2MTViOGIiLCJhcHBJZCI6MjMyLCJ0b2tlbklkIjoiMWQzMzk2NzQtNjg5MC00ZjlmLTk4ZjYtMTAxYWYwZjI4NmMxIiwiaWF0IjoxNzE1NjM3MDc2fQ.l0P4OAWQvDMOjebNE8xRBNnU8nNsgiJQdyFH5KEAmPk"
$ws.Range("F3").Value = "mp4n4f2mq3ca"
$ws.Range("G3").Value = "akhil"
$ws.Range("H3").Value = "48298ea0-240a-11ee-9c04-b37209615b8b"
# Note: I3 (LichessReplaced) intentionally left blank - it was not set in the edit.
